$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing sheet to "trialwise" and set its selection ---
$trial = $wb.Worksheets.Item(1)
$trial.Name = "trialwise"
$trial.Range("B39").Select()

# --- Step 2: insert the new "blockwise" sheet (Excel inserts it before the active sheet) ---
$block = $wb.Worksheets.Add()
$block.Name = "blockwise"

# --- Step 3: populate the blockwise summary table (rows 2-17 first, header row last, so
#             the shared-string table picks up "Emotion"/"Why"/"Action"/"How" before
#             "Stimulus"/"Question" -- matching the order new strings were appended) ---
$rows = @(
    @("Is the person admiring someone?", "admiring?", "Emotion", "Why"),
    @("Is the person expressing self-doubt?", "self-doubt?", "Emotion", "Why"),
    @("Is the person in an argument?", "argument?", "Emotion", "Why"),
    @("Is the person proud of themselves?", "proud?", "Emotion", "Why"),
    @("Is the person competing against others?", "competing?", "Action", "Why"),
    @("Is the person concerned with their health?", "healthy?", "Action", "Why"),
    @("Is the person helping someone?", "helping?", "Action", "Why"),
    @("Is the person protecting themselves?", "self-protection?", "Action", "Why"),
    @("Is the person looking at the camera?", "looking at camera?", "Emotion", "How"),
    @("Is the person looking to their side?", "looking to side?", "Emotion", "How"),
    @("Is the person opening their mouth?", "open mouth?", "Emotion", "How"),
    @("Is the person smiling?", "smiling?", "Emotion", "How"),
    @("Is the person lifting something?", "lifting?", "Action", "How"),
    @("Is the person pressing a button?", "pressing button?", "Action", "How"),
    @("Is the person reaching for something?", "reaching?", "Action", "How"),
    @("Is the person using both hands?", "both hands?", "Action", "How")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $vals = $rows[$i]
    $block.Cells.Item($r, 1).Value = $vals[0]
    $block.Cells.Item($r, 2).Value = $vals[1]
    $block.Cells.Item($r, 3).Value = $vals[2]
    $block.Cells.Item($r, 4).Value = $vals[3]
}

# Header row written last so "Stimulus"/"Question" are the final two new shared strings.
$block.Cells.Item(1, 1).Value = "Question Cue"
$block.Cells.Item(1, 2).Value = "Reminder Cue"
$block.Cells.Item(1, 3).Value = "Stimulus"
$block.Cells.Item(1, 4).Value = "Question"

# --- Step 4: formatting ---
# Body rows: column A left-aligned Arial 11, columns B-D centered Arial 11 (same look as trialwise)
$bodyRange = $block.Range("A2:D17")
$bodyRange.Font.Name = "Arial"
$bodyRange.Font.Size = 11
$bodyRange.Font.Bold = $false
$bodyRange.HorizontalAlignment = -4108

$block.Range("A2:A17").HorizontalAlignment = -4131
$block.Range("B2:D17").HorizontalAlignment = -4108

# Header row: bold Arial, shaded fill
$headerA = $block.Range("A1")
$headerA.Font.Name = "Arial"
$headerA.Font.Size = 12
$headerA.Font.Bold = $false
$headerA.HorizontalAlignment = -4131

$headerB = $block.Range("B1")
$headerB.Font.Name = "Arial"
$headerB.Font.Size = 12
$headerB.Font.Bold = $false
$headerB.HorizontalAlignment = -4108

$headerCD = $block.Range("C1:D1")
$headerCD.Font.Name = "Arial"
$headerCD.Font.Bold = $true
$headerCD.Font.Size = 11
$headerCD.Font.Color = 0
$headerCD.Interior.Color = 14277081
$headerCD.HorizontalAlignment = -4108

# --- Step 5: column widths ---
$block.Columns.Item(1).ColumnWidth = 35.33203125
$block.Columns.Item(2).ColumnWidth = 16.83203125

# --- Step 6: final selection / active cell on blockwise ---
$block.Range("G13").Select()
